$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Diary update: append two new rows (14 and 15) for 17/05/2024 (serial 45429) ---
# Row 14: Telecomunicazioni 1
$ws.Cells.Item(14, 1).Value = 45429
$ws.Cells.Item(14, 2).Value = "Telecomunicazioni 1"
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 1

# Row 15: Telecomunicazioni 2
$ws.Cells.Item(15, 1).Value = 45429
$ws.Cells.Item(15, 2).Value = "Telecomunicazioni 2"
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 1

# Match the look (number format / borders / font) of the row above (13) for the new rows
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 14/15 use the same compact row height as the other data rows (13.8pt)
$ws.Rows.Item(14).RowHeight = 13.8
$ws.Rows.Item(15).RowHeight = 13.8

# --- Column width tweaks ---
# Wider default column width for the sheet
$ws.StandardWidth = 19.2890625

# Column A a bit narrower (target stored width ~12.24), columns C:F slightly
# narrower too (target stored width ~12.46) - closest achievable char widths
$ws.Columns.Item(1).ColumnWidth = 11.25
$ws.Range("C1:F1").ColumnWidth = 11.584

# --- Update the active cell / selection ---
$ws.Range("H11").Select()
